$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$clothingData = @(
    @{Row=2; Value="Halter,Blazer"},
    @{Row=3; Value="Parka,Caftan"},
    @{Row=4; Value="Jumpsuit,Halter"},
    @{Row=5; Value="Blouse,Jumpsuit"},
    @{Row=6; Value="Jumpsuit,Halter"},
    @{Row=7; Value="Trunks,Caftan"},
    @{Row=8; Value="Parka,Caftan"},
    @{Row=9; Value="Blazer,Top"},
    @{Row=10; Value="Gauchos,Parka"},
    @{Row=11; Value="Jumpsuit,Trunks"},
    @{Row=12; Value="Jumpsuit,Dress"},
    @{Row=13; Value="Parka,Blouse"},
    @{Row=14; Value="Jumpsuit,Blouse"},
    @{Row=15; Value="Caftan,Jumpsuit"},
    @{Row=16; Value="Jumpsuit,Blouse"},
    @{Row=17; Value="Parka,Jumpsuit"},
    @{Row=18; Value="Caftan,Parka"},
    @{Row=19; Value="Jumpsuit,Blouse"},
    @{Row=20; Value="Jumpsuit,Tee"},
    @{Row=21; Value="Parka,Caftan"},
    @{Row=22; Value="Jumpsuit,Dress"},
    @{Row=23; Value="Halter,Blouse"},
    @{Row=24; Value="Jumpsuit,Dress"},
    @{Row=25; Value="Parka,Caftan"},
    @{Row=26; Value="Parka,Caftan"},
    @{Row=27; Value="Parka,Caftan"},
    @{Row=28; Value="Halter,Blazer"},
    @{Row=29; Value="Blazer,Halter"},
    @{Row=30; Value="Jumpsuit,Kaftan"},
    @{Row=31; Value="Jumpsuit,Parka"},
    @{Row=32; Value="Halter,Caftan"},
    @{Row=33; Value="Halter,Blouse"},
    @{Row=34; Value="Blouse,Jumpsuit"},
    @{Row=35; Value="Halter,Jumpsuit"},
    @{Row=36; Value="Dress,Jumpsuit"},
    @{Row=37; Value="Kaftan,Jumpsuit"},
    @{Row=38; Value="Jumpsuit,Sweatpants"},
    @{Row=39; Value="Jumpsuit,Kaftan"},
    @{Row=40; Value="Jumpsuit,Blouse"},
    @{Row=41; Value="Jumpsuit,Kaftan"},
    @{Row=42; Value="Jumpsuit,Blouse"},
    @{Row=43; Value="Jumpsuit,Blouse"},
    @{Row=44; Value="Jumpsuit,Blouse"},
    @{Row=45; Value="Parka,Blouse"},
    @{Row=46; Value="Halter,Blouse"},
    @{Row=47; Value="Blouse,Halter"},
    @{Row=48; Value="Jumpsuit,Kaftan"},
    @{Row=49; Value="Parka,Blouse"},
    @{Row=50; Value="Halter,Parka"},
    @{Row=51; Value="Blouse,Jumpsuit"},
    @{Row=52; Value="Parka,Halter"},
    @{Row=53; Value="Parka,Blouse"},
    @{Row=54; Value="Parka,Jumpsuit"},
    @{Row=55; Value="Jumpsuit,Kaftan"},
    @{Row=56; Value="Jumpsuit,Kaftan"}
)

foreach ($entry in $clothingData) {
    $ws.Cells.Item($entry.Row, 7).Value = $entry.Value
}

$wb.Save()
